$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: ", mas ao todos temos .......... tipos de girassóis."
#      -> ", mas ao todos temos cerca de 70 espécies de girassóis."
# split across three runs (same Arial/auto/24/24 formatting as the original
# run) matching the way the author retyped the sentence.
# ---------------------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.Execute(
    ", mas ao todos temos .......... tipos de girassóis.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "",
    0
)
$start = $findRng.Start

# Remove the old sentence entirely.
$findRng.Text = ""

$t1 = ", mas ao todos temos "
$t2 = "cerca de 70 espécies "
$t3 = "de girassóis."

$p0 = $start
$p1 = $p0 + $t1.Length
$p2 = $p1 + $t2.Length

# Run 1 - merges into the (unchanged) preceding run, which already carries
# the same run formatting, exactly as in the target document.
$r1 = $d.Range($p0, $p0)
$r1.Text = $t1

# Run 2
$r2 = $d.Range($p1, $p1)
$r2.Text = $t2

# Run 3
$r3 = $d.Range($p2, $p2)
$r3.Text = $t3

# Force a hard run boundary between run1/run2 (both otherwise share identical
# formatting and would be silently re-coalesced) by toggling and restoring a
# character attribute over exactly the run-2 span.
$boundary = $d.Range($p1, $p2)
$boundary.Font.Bold = 1
$boundary.Font.Bold = 0

# ---------------------------------------------------------------------------
# Edit 2: drop the "SemEspaamento" (No Spacing) character style from the
# single-space run that follows "explique os impactos."
# ---------------------------------------------------------------------------
$impactRng = $d.Content
$impactRng.Find.Execute(
    "explique os impactos.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "",
    0
)
$impactEnd = $impactRng.End
$spaceRng = $d.Range($impactEnd, $impactEnd + 1)

# Re-create the run from scratch so it picks up a plain rPr with no rStyle,
# then force it to stay a distinct run (it would otherwise re-merge with the
# identically-formatted "explique os impactos." run).
$spaceRng.Delete()
$insertRng = $d.Range($impactEnd, $impactEnd)
$insertRng.InsertAfter(" ")

$newSpaceRng = $d.Range($impactEnd, $impactEnd + 1)
$newSpaceRng.Font.Bold = 1
$newSpaceRng.Font.Bold = 0
